# feat: add 2022-Q1 data
#
# Inserts a new worksheet "2022-Q1" (with the per-fund holding breakdown,
# same column layout as the existing "2021-Q1" sheet) positioned between
# "2021-Q1" and "总计", and records the new quarter's summary row at the
# top of the "总计" sheet's data table.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new "2022-Q1" worksheet right after "2021-Q1".
# ---------------------------------------------------------------------
$q1Sheet = $wb.Worksheets.Item(1)
$newSheet = $wb.Worksheets.Add($null, $q1Sheet)
$newSheet.Name = "2022-Q1"

# Pick up the bold/centered/bordered header & index-column look used by
# the "总计" sheet (copy formats only, values are filled in afterwards).
$wb.Worksheets.Item("总计").Range("B1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$wb.Worksheets.Item("总计").Range("A2").Copy()
$newSheet.Range("A2:A3").PasteSpecial(-4122)

# Header row
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Columns B, D, E, F, G hold text that looks numeric (fund codes with
# leading zeros, decimal figures stored as strings in the source data) -
# force text format first so COM doesn't silently coerce them to numbers.
# (Multi-area ranges only honour the format on the first area, so set
# each contiguous block separately.)
$newSheet.Range("B2:B3").NumberFormat = "@"
$newSheet.Range("D2:G3").NumberFormat = "@"

# Row 2 - fund A share
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "009837"
$newSheet.Range("C2").Value = "华夏磐锐一年定期开放混合A"
$newSheet.Range("D2").Value = "16.45"
$newSheet.Range("E2").Value = "79.44"
$newSheet.Range("F2").Value = "2.68"
$newSheet.Range("G2").Value = "0.4409"
$newSheet.Range("H2").Value = 7

# Row 3 - fund C share
$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "009838"
$newSheet.Range("C3").Value = "华夏磐锐一年定期开放混合C"
$newSheet.Range("D3").Value = "0.44"
$newSheet.Range("E3").Value = "79.44"
$newSheet.Range("F3").Value = "2.68"
$newSheet.Range("G3").Value = "0.0118"
$newSheet.Range("H3").Value = 7

# ---------------------------------------------------------------------
# 2. Add the 2022-Q1 summary row to the "总计" sheet (newest quarter on
#    top), pushing the existing 2021-Q1 row down.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()

# The inserted row inherits formatting from the row above (the header);
# the data cells B2:D2 should be plain, like the rest of the table, while
# A2 should keep the bold/centered/bordered index-column look - copy that
# from A3 (the original row, now shifted down), which still has it.
$totalSheet.Range("B2:D2").ClearFormats()
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 2
$totalSheet.Range("D2").Value = 0.45

# Renumber the pre-existing 2021-Q1 row, now shifted to row 3.
$totalSheet.Range("A3").Value = 1
